$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix ("bug in KVIK dyr samlet"): the "kontrol" scenario rows were
# emitted with no data (blank TotGoednabDyr_kt_aar) and should not be in
# the output table at all. Remove both "kontrol" rows (svin + kvaeg).
$ws.Range("A8:C9").EntireRow.Delete()
